$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.801.98"
$ws.Range("E2").Value = "  +2.66%  "

$ws.Range("D3").Value = "3.837.93"
$ws.Range("E3").Value = "  +1.30%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.35%  "

$ws.Range("D5").Value = "'635.23"
$ws.Range("E5").Value = "  +5.79%  "

$ws.Range("D6").Value = "'166.79"
$ws.Range("E6").Value = "  +0.96%  "

$ws.Range("D7").Value = "3.835.93"
$ws.Range("E7").Value = "  +1.33%  "

$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").Value = "'0.522"
$ws.Range("E9").Value = "  +0.88%  "

$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = "  +2.32%  "

$ws.Range("D11").Value = "'0.455"
$ws.Range("E11").Value = "  +0.81%  "

$ws.Range("D12").Value = "'6.67"
$ws.Range("E12").Value = "  +2.97%  "

$ws.Range("E13").Value = "  +1.19%  "

$ws.Range("D14").Value = "'36.08"
$ws.Range("E14").Value = "  +1.02%  "

$ws.Range("D15").Value = "4.477.11"

$ws.Range("D16").Value = "3.867.23"
$ws.Range("E16").Value = "  +1.61%  "

$ws.Range("D17").Value = "69.681.48"
$ws.Range("E17").Value = "  +2.41%  "

$ws.Range("E18").Value = "  -1.39%  "

$ws.Range("E19").Value = "  +1.19%  "

$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").Value = "'468.62"
$ws.Range("E21").Value = "  +1.68%  "

$ws.Range("D22").Value = "'9.75"
$ws.Range("E22").Value = "  +0.58%  "

$ws.Range("D23").Value = "'0.710"
$ws.Range("E23").Value = "  +1.71%  "

$ws.Range("E24").Value = "  +2.45%  "

$ws.Range("D25").Value = "'83.78"
$ws.Range("E25").Value = "  +1.22%  "

$ws.Range("D26").Value = "'2.19"
$ws.Range("E26").Value = "  +3.90%  "

$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("D28").Value = "'10.11"
$ws.Range("E28").Value = "  +1.29%  "

$ws.Range("D30").Value = "3.983.79"

$ws.Range("E31").Value = "  +2.87%  "

$ws.Range("D32").Value = "'7.33"
$ws.Range("E32").Value = "  -0.70%  "

$ws.Range("E33").Value = "  +0.30%  "

$ws.Range("D34").Value = "'29.33"
$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("D35").Value = "3.777.95"

$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").Value = "'9.09"
$ws.Range("E37").Value = "  +1.10%  "

$ws.Range("E38").Value = "  +3.07%  "

$ws.Range("E39").Value = "  +8.30%  "

$ws.Range("D40").Value = "'3.51"
$ws.Range("E40").Value = "  +7.80%  "

$ws.Range("D41").Value = "'5.93"
$ws.Range("E41").Value = "  +2.79%  "

$ws.Range("D42").Value = "'0.983"
$ws.Range("E42").Value = "  -0.60%  "

$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.12%  "

$ws.Range("D45").Value = "'157.49"
$ws.Range("E45").Value = "  +3.91%  "

$ws.Range("D46").Value = "'44.13"
$ws.Range("E46").Value = "  +2.49%  "

$ws.Range("E47").Value = "  +1.14%  "

$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'47.32"
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.94"
$ws.Range("E49").Value = "  +3.64%  "

$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "'1.41"
$ws.Range("E50").Value = "  +5.30%  "

$ws.Range("D51").Value = "'8.47"
$ws.Range("E51").Value = "  +1.41%  "
